$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# @@ -3730,25 +3730,25 @@  sheet=ALC
$ws_ALC.Range("H62").Value = 5463.75
$ws_ALC.Range("I62").Value = 1677.5
$ws_ALC.Range("J62").Value = 9250
$ws_ALC.Range("K62").Value = 1677.5
$ws_ALC.Range("L62").Value = 9250
$ws_ALC.Range("M62").Value = -1053.5
$ws_ALC.Range("N62").Value = -10498

# @@ -3883,25 +3883,25 @@  sheet=ALC
$ws_ALC.Range("H65").Value = 5463.75
$ws_ALC.Range("I65").Value = 1677.5
$ws_ALC.Range("J65").Value = 9250
$ws_ALC.Range("K65").Value = 8387.5
$ws_ALC.Range("L65").Value = 46250
$ws_ALC.Range("M65").Value = -5267.5
$ws_ALC.Range("N65").Value = -52490

# @@ -4942,25 +4942,25 @@  sheet=ALC
$ws_ALC.Range("H86").Value = 2003.5238
$ws_ALC.Range("I86").Value = 1790.7693
$ws_ALC.Range("J86").Value = 2349.25
$ws_ALC.Range("K86").Value = 1790.7693
$ws_ALC.Range("L86").Value = 2349.25
$ws_ALC.Range("M86").Value = -667.7692999999999
$ws_ALC.Range("N86").Value = -4595.25

# @@ -4994,22 +4994,22 @@  sheet=ALC
$ws_ALC.Range("H87").Value = 38766
$ws_ALC.Range("J87").Value = 38766
$ws_ALC.Range("L87").Value = 38766
$ws_ALC.Range("N87").Value = -41262

# @@ -5092,25 +5092,25 @@  sheet=ALC
$ws_ALC.Range("H89").Value = 2003.5238
$ws_ALC.Range("I89").Value = 1790.7693
$ws_ALC.Range("J89").Value = 2349.25
$ws_ALC.Range("K89").Value = 8953.8465
$ws_ALC.Range("L89").Value = 11746.25
$ws_ALC.Range("M89").Value = -3337.8465
$ws_ALC.Range("N89").Value = -22978.25

# @@ -5144,22 +5144,22 @@  sheet=ALC
$ws_ALC.Range("H90").Value = 38766
$ws_ALC.Range("J90").Value = 38766
$ws_ALC.Range("L90").Value = 116298
$ws_ALC.Range("N90").Value = -128778

# @@ -5545,25 +5545,25 @@  sheet=ALC
$ws_ALC.Range("H98").Value = 1313.5333
$ws_ALC.Range("I98").Value = 895.4583
$ws_ALC.Range("J98").Value = 2985.8333
$ws_ALC.Range("K98").Value = 895.4583
$ws_ALC.Range("L98").Value = 2985.8333
$ws_ALC.Range("M98").Value = 602.5417
$ws_ALC.Range("N98").Value = -5981.8333

# @@ -6745,25 +6745,25 @@  sheet=ALC
$ws_ALC.Range("H122").Value = 1313.5333
$ws_ALC.Range("I122").Value = 895.4583
$ws_ALC.Range("J122").Value = 2985.8333
$ws_ALC.Range("K122").Value = 2686.3749
$ws_ALC.Range("L122").Value = 8957.499899999999
$ws_ALC.Range("M122").Value = -236.3748999999998
$ws_ALC.Range("N122").Value = -13857.4999

# @@ -7097,7 +7097,7 @@  sheet=ALC
$ws_ALC.Range("H129").Value = 711.96

# @@ -7501,22 +7501,22 @@  sheet=ALC
$ws_ALC.Range("H137").Value = 30232.922
$ws_ALC.Range("I137").Value = 1091
$ws_ALC.Range("K137").Value = 3273
$ws_ALC.Range("M137").Value = -723

# @@ -9325,22 +9325,22 @@  sheet=ARM
$ws_ARM.Range("H32").Value = 21292474
$ws_ARM.Range("I32").Value = 29416720
$ws_ARM.Range("K32").Value = 29416720
$ws_ARM.Range("M32").Value = -29416433

# @@ -13111,25 +13111,25 @@  sheet=ARM
$ws_ARM.Range("H110").Value = 4300
$ws_ARM.Range("I110").Value = 3477.7778
$ws_ARM.Range("J110").Value = 5357.143
$ws_ARM.Range("K110").Value = 3477.7778
$ws_ARM.Range("L110").Value = 5357.143
$ws_ARM.Range("M110").Value = -1432.7778
$ws_ARM.Range("N110").Value = -9447.143

# @@ -13696,22 +13696,22 @@  sheet=ARM
$ws_ARM.Range("H122").Value = 980.25
$ws_ARM.Range("I122").Value = 887.5454999999999
$ws_ARM.Range("K122").Value = 2662.6365
$ws_ARM.Range("M122").Value = -212.6364999999996

# @@ -13748,22 +13748,22 @@  sheet=ARM
$ws_ARM.Range("H123").Value = 20500
$ws_ARM.Range("J123").Value = 20500
$ws_ARM.Range("L123").Value = 20500
$ws_ARM.Range("N123").Value = -30300

# @@ -14189,25 +14189,25 @@  sheet=ARM
$ws_ARM.Range("H132").Value = 1754.4412
$ws_ARM.Range("I132").Value = 1173.88
$ws_ARM.Range("J132").Value = 3367.111
$ws_ARM.Range("K132").Value = 3521.64
$ws_ARM.Range("L132").Value = 10101.333
$ws_ARM.Range("M132").Value = -991.6400000000003
$ws_ARM.Range("N132").Value = -15161.333

# @@ -18895,25 +18895,25 @@  sheet=BSM
$ws_BSM.Range("H86").Value = 1878.3
$ws_BSM.Range("I86").Value = 1833
$ws_BSM.Range("J86").Value = 1984
$ws_BSM.Range("K86").Value = 1833
$ws_BSM.Range("L86").Value = 1984
$ws_BSM.Range("M86").Value = -710
$ws_BSM.Range("N86").Value = -4230

# @@ -19039,25 +19039,25 @@  sheet=BSM
$ws_BSM.Range("H89").Value = 1878.3
$ws_BSM.Range("I89").Value = 1833
$ws_BSM.Range("J89").Value = 1984
$ws_BSM.Range("K89").Value = 9165
$ws_BSM.Range("L89").Value = 9920
$ws_BSM.Range("M89").Value = -3549
$ws_BSM.Range("N89").Value = -21152

# @@ -19826,25 +19826,25 @@  sheet=BSM
$ws_BSM.Range("H105").Value = 3014.743
$ws_BSM.Range("I105").Value = 2010
$ws_BSM.Range("J105").Value = 3029.3044
$ws_BSM.Range("K105").Value = 2010
$ws_BSM.Range("L105").Value = 3029.3044
$ws_BSM.Range("M105").Value = -263
$ws_BSM.Range("N105").Value = -6523.3044

# @@ -19924,25 +19924,22 @@  sheet=BSM
$ws_BSM.Range("H107").Value = 3663.4167
$ws_BSM.Range("I107").Value = 3663.4167
$ws_BSM.Range("J107").Value = 0
$ws_BSM.Range("K107").Value = 3663.4167
$ws_BSM.Range("L107").Value = 0
$ws_BSM.Range("M107").Value = -1743.4167
$ws_BSM.Range("N107").ClearContents()

# @@ -29183,25 +29180,25 @@  sheet=CUL
$ws_CUL.Range("H12").Value = 6.142857
$ws_CUL.Range("I12").Value = 13.75
$ws_CUL.Range("J12").Value = 3.1
$ws_CUL.Range("K12").Value = 41.25
$ws_CUL.Range("L12").Value = 9.300000000000001
$ws_CUL.Range("M12").Value = 131.75
$ws_CUL.Range("N12").Value = -355.3

# @@ -30487,25 +30484,25 @@  sheet=CUL
$ws_CUL.Range("H38").Value = 55.285713
$ws_CUL.Range("I38").Value = 30.333334
$ws_CUL.Range("J38").Value = 74
$ws_CUL.Range("K38").Value = 91.00000199999999
$ws_CUL.Range("L38").Value = 222
$ws_CUL.Range("M38").Value = 255.999998
$ws_CUL.Range("N38").Value = -916

# @@ -35246,25 +35243,25 @@  sheet=CUL
$ws_CUL.Range("H132").Value = 22096520
$ws_CUL.Range("I132").Value = 694.6667
$ws_CUL.Range("J132").Value = 35354016
$ws_CUL.Range("K132").Value = 6252.0003
$ws_CUL.Range("L132").Value = 318186144
$ws_CUL.Range("M132").Value = -3722.0003
$ws_CUL.Range("N132").Value = -318191204

# @@ -35454,22 +35451,22 @@  sheet=CUL
$ws_CUL.Range("H136").Value = 4402
$ws_CUL.Range("I136").Value = 3114.8333
$ws_CUL.Range("K136").Value = 9344.499899999999
$ws_CUL.Range("M136").Value = -4244.499899999999

# @@ -35610,25 +35607,25 @@  sheet=CUL
$ws_CUL.Range("H139").Value = 2802.0908
$ws_CUL.Range("I139").Value = 1140.1111
$ws_CUL.Range("J139").Value = 3952.6924
$ws_CUL.Range("K139").Value = 3420.3333
$ws_CUL.Range("L139").Value = 11858.0772
$ws_CUL.Range("M139").Value = 1719.6667
$ws_CUL.Range("N139").Value = -22138.0772

# @@ -39174,25 +39171,25 @@  sheet=GSM
$ws_GSM.Range("H70").Value = 5428.5776
$ws_GSM.Range("I70").Value = 5457.5527
$ws_GSM.Range("J70").Value = 5271.2856
$ws_GSM.Range("K70").Value = 5457.5527
$ws_GSM.Range("L70").Value = 5271.2856
$ws_GSM.Range("M70").Value = -5187.5527
$ws_GSM.Range("N70").Value = -5811.2856

# @@ -39321,25 +39318,25 @@  sheet=GSM
$ws_GSM.Range("H73").Value = 5428.5776
$ws_GSM.Range("I73").Value = 5457.5527
$ws_GSM.Range("J73").Value = 5271.2856
$ws_GSM.Range("K73").Value = 5457.5527
$ws_GSM.Range("L73").Value = 5271.2856
$ws_GSM.Range("M73").Value = -4521.5527
$ws_GSM.Range("N73").Value = -7143.2856

# @@ -39667,22 +39664,22 @@  sheet=GSM
$ws_GSM.Range("H80").Value = 3167.2
$ws_GSM.Range("I80").Value = 2694.6316
$ws_GSM.Range("K80").Value = 2694.6316
$ws_GSM.Range("M80").Value = -1696.6316

# @@ -39817,22 +39814,22 @@  sheet=GSM
$ws_GSM.Range("H83").Value = 3167.2
$ws_GSM.Range("I83").Value = 2694.6316
$ws_GSM.Range("K83").Value = 13473.158
$ws_GSM.Range("M83").Value = -8481.158000000001

# @@ -43455,22 +43452,25 @@  sheet=LTW
$ws_LTW.Range("H16").Value = 770.9375
$ws_LTW.Range("I16").Value = 666.7857
$ws_LTW.Range("J16").Value = 1500
$ws_LTW.Range("K16").Value = 666.7857
$ws_LTW.Range("L16").Value = 1500
$ws_LTW.Range("M16").Value = -496.7857
$ws_LTW.Range("N16").Value = -1840

# @@ -46907,22 +46907,22 @@  sheet=LTW
$ws_LTW.Range("H87").Value = 33500
$ws_LTW.Range("J87").Value = 33500
$ws_LTW.Range("L87").Value = 33500
$ws_LTW.Range("N87").Value = -35746

# @@ -46956,25 +46956,25 @@  sheet=LTW
$ws_LTW.Range("H88").Value = 21482
$ws_LTW.Range("I88").Value = 14428
$ws_LTW.Range("J88").Value = 23833.334
$ws_LTW.Range("K88").Value = 14428
$ws_LTW.Range("L88").Value = 23833.334
$ws_LTW.Range("M88").Value = -14000
$ws_LTW.Range("N88").Value = -24689.334

# @@ -47057,22 +47057,22 @@  sheet=LTW
$ws_LTW.Range("H90").Value = 33500
$ws_LTW.Range("J90").Value = 33500
$ws_LTW.Range("L90").Value = 100500
$ws_LTW.Range("N90").Value = -111732

# @@ -47106,25 +47106,25 @@  sheet=LTW
$ws_LTW.Range("H91").Value = 21482
$ws_LTW.Range("I91").Value = 14428
$ws_LTW.Range("J91").Value = 23833.334
$ws_LTW.Range("K91").Value = 14428
$ws_LTW.Range("L91").Value = 23833.334
$ws_LTW.Range("M91").Value = -12946
$ws_LTW.Range("N91").Value = -26797.334

# @@ -48086,22 +48086,19 @@  sheet=LTW
$ws_LTW.Range("H111").Value = 0
$ws_LTW.Range("J111").Value = 0
$ws_LTW.Range("L111").Value = 0
$ws_LTW.Range("N111").ClearContents()

# @@ -49794,19 +49791,25 @@  sheet=WVR
$ws_WVR.Range("H4").Value = 352.75
$ws_WVR.Range("I4").Value = 202
$ws_WVR.Range("J4").Value = 403
$ws_WVR.Range("K4").Value = 202
$ws_WVR.Range("L4").Value = 403
$ws_WVR.Range("M4").Value = -89
$ws_WVR.Range("N4").Value = -629

# @@ -50887,25 +50890,22 @@  sheet=WVR
$ws_WVR.Range("H26").Value = 2500333.2
$ws_WVR.Range("I26").Value = 0
$ws_WVR.Range("J26").Value = 2500333.2
$ws_WVR.Range("K26").Value = 0
$ws_WVR.Range("L26").Value = 2500333.2
$ws_WVR.Range("M26").ClearContents()
$ws_WVR.Range("N26").Value = -2500919.2

# @@ -52633,25 +52633,25 @@  sheet=WVR
$ws_WVR.Range("H62").Value = 3625
$ws_WVR.Range("I62").Value = 2500
$ws_WVR.Range("J62").Value = 4000
$ws_WVR.Range("K62").Value = 2500
$ws_WVR.Range("L62").Value = 4000
$ws_WVR.Range("M62").Value = -1876
$ws_WVR.Range("N62").Value = -5248

# @@ -52780,25 +52780,25 @@  sheet=WVR
$ws_WVR.Range("H65").Value = 3625
$ws_WVR.Range("I65").Value = 2500
$ws_WVR.Range("J65").Value = 4000
$ws_WVR.Range("K65").Value = 12500
$ws_WVR.Range("L65").Value = 20000
$ws_WVR.Range("M65").Value = -9380
$ws_WVR.Range("N65").Value = -26240

# @@ -55108,25 +55108,25 @@  sheet=WVR
$ws_WVR.Range("H113").Value = 287.65216
$ws_WVR.Range("I113").Value = 292.65
$ws_WVR.Range("J113").Value = 254.33333
$ws_WVR.Range("K113").Value = 877.9499999999999
$ws_WVR.Range("L113").Value = 762.99999
$ws_WVR.Range("M113").Value = 1292.05
$ws_WVR.Range("N113").Value = -5102.99999

Write-Output "Applied all Shinryu_Profits updates"
